# Insert a new weekly record row at row 28, pushing all subsequent
# rows down by one (this updates the sheet's dimension automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("28:28").Insert()

$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 45177
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 100112035
$ws.Cells.Item(28, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 52
$ws.Cells.Item(28, 11).Value = 19000
$ws.Cells.Item(28, 12).Value = 21000
$ws.Cells.Item(28, 13).Value = 20000
$ws.Cells.Item(28, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(28, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(28, 16).Value = 1333
$ws.Cells.Item(28, 17).Value = 15
$ws.Cells.Item(28, 18).Value = "Hortaliza"
